$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newRow = 78

# Force column A's new cell to be stored as text (not auto-parsed as a date)
# by giving it a Text number format before assigning the literal value.
$ws.Cells.Item($newRow, 1).NumberFormat = "@"
$ws.Cells.Item($newRow, 1).Value = "03-12-2025"
$ws.Cells.Item($newRow, 2).Value = "The price of gold in India today is ₹13,058 per gram for 24 karat gold, ₹11,970 per gram for 22 karat gold and ₹9,794 per gram for 18 karat gold (also called 999 gold)."

# Re-apply the same formatting (borders/wrap) used by the rest of the table
# so the new row visually matches the previous data rows.
$ws.Range("A77:B77").Copy()
$ws.Range("A78:B78").PasteSpecial(-4122)
